$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-47 down to 43-48
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new price record
$ws.Cells.Item(42, 1).Value = 1
$ws.Cells.Item(42, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(42, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(42, 4).Value = 44491
$ws.Cells.Item(42, 5).Value = 15
$ws.Cells.Item(42, 6).Value = 100112021
$ws.Cells.Item(42, 7).Value = "Ají"
$ws.Cells.Item(42, 8).Value = "Inferno"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 120
$ws.Cells.Item(42, 11).Value = 37000
$ws.Cells.Item(42, 12).Value = 38000
$ws.Cells.Item(42, 13).Value = 37500
$ws.Cells.Item(42, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(42, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(42, 16).Value = 2500
$ws.Cells.Item(42, 17).Value = 15
$ws.Cells.Item(42, 18).Value = "Hortaliza"
